$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextCell "D2" "28.467.73"
Set-TextCell "E2" "  -3.64%  "
Set-TextCell "D3" "1.957.25"
Set-TextCell "E3" "  -2.20%  "
Set-TextCell "E4" "  -0.44%  "
Set-TextCell "D5" "321.50"
Set-TextCell "E5" "  -2.40%  "
Set-TextCell "E6" "  -0.23%  "
Set-TextCell "D7" "0.4757"
Set-TextCell "E7" "  -4.96%  "
Set-TextCell "D8" "0.4058"
Set-TextCell "E8" "  -3.96%  "
Set-TextCell "D9" "53.17"
Set-TextCell "E9" "  -2.15%  "
Set-TextCell "D10" "0.08434"
Set-TextCell "E10" "  -6.04%  "
Set-TextCell "D11" "1.058"
Set-TextCell "E11" "  -5.06%  "
Set-TextCell "D12" "22.14"
Set-TextCell "E12" "  -4.85%  "
Set-TextCell "D13" "1.981.91"
Set-TextCell "E13" "  -3.66%  "
Set-TextCell "D14" "7.609"
Set-TextCell "E14" "  -5.45%  "
Set-TextCell "D15" "6.169"
Set-TextCell "E15" "  -4.27%  "
Set-TextCell "D16" "1.012"
Set-TextCell "E16" "  -0.19%  "
Set-TextCell "B17" "Litecoin"
Set-TextCell "C17" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D17" "89.35"
Set-TextCell "E17" "  -5.13%  "
Set-TextCell "B18" "ShibaInu"
Set-TextCell "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D18" "0.00001073"
Set-TextCell "E18" "  -3.41%  "
Set-TextCell "D19" "0.06611"
Set-TextCell "E19" "  -1.16%  "
Set-TextCell "D20" "18.69"
Set-TextCell "E20" "  -4.69%  "
Set-TextCell "E21" "  -0.20%  "
Set-TextCell "D22" "5.818"
Set-TextCell "E22" "  -2.25%  "
Set-TextCell "D23" "28.501.35"
Set-TextCell "E23" "  -3.67%  "
Set-TextCell "D24" "11.57"
Set-TextCell "E24" "  -3.30%  "
Set-TextCell "D25" "2.292"
Set-TextCell "E25" "  -0.49%  "
Set-TextCell "D26" "2.208.00"
Set-TextCell "E26" "  -2.51%  "
Set-TextCell "D27" "154.90"
Set-TextCell "E27" "  -1.99%  "
Set-TextCell "D28" "20.20"
Set-TextCell "E28" "  -2.52%  "
Set-TextCell "D29" "5.953"
Set-TextCell "E29" "  -5.88%  "
Set-TextCell "D30" "2.158"
Set-TextCell "E30" "  -5.97%  "
Set-TextCell "D31" "123.49"
Set-TextCell "E31" "  -3.45%  "
Set-TextCell "D32" "0.9818"
Set-TextCell "E32" "  -6.79%  "
Set-TextCell "D33" "0.09598"
Set-TextCell "E33" "  -3.41%  "
Set-TextCell "D34" "1.446"
Set-TextCell "E34" "  -6.50%  "
Set-TextCell "B35" "HuobiToken"
Set-TextCell "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D35" "3.662"
Set-TextCell "E35" "  -3.59%  "
Set-TextCell "B36" "Filecoin"
Set-TextCell "C36" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D36" "5.581"
Set-TextCell "E36" "  -4.15%  "
Set-TextCell "D37" "0.02336"
Set-TextCell "E37" "  -4.70%  "
Set-TextCell "D38" "8.802"
Set-TextCell "E38" "  -5.03%  "
Set-TextCell "D39" "0.06220"
Set-TextCell "E39" "  -2.76%  "
Set-TextCell "D40" "1.258"
Set-TextCell "E40" "  -3.37%  "
Set-TextCell "D41" "0.6217"
Set-TextCell "E41" "  -5.02%  "
Set-TextCell "D42" "11.13"
Set-TextCell "E42" "  -4.20%  "
Set-TextCell "E43" "  -0.28%  "
Set-TextCell "D44" "0.1920"
Set-TextCell "D45" "1.337"
Set-TextCell "E45" "  +2.85%  "
Set-TextCell "D46" "0.5962"
Set-TextCell "E46" "  -5.86%  "
Set-TextCell "D47" "13.05"
Set-TextCell "E47" "  -3.14%  "
Set-TextCell "D48" "2.053"
Set-TextCell "E48" "  -6.21%  "
Set-TextCell "D49" "3.400"
Set-TextCell "E49" "  -3.03%  "
Set-TextCell "D50" "0.00000000327"
Set-TextCell "E50" "  -3.58%  "
Set-TextCell "D51" "0.06829"
Set-TextCell "E51" "  -2.20%  "
